$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1): rename some headers, add 4 new trailing headers
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 7).Value  = "점수(룰)"
$ws.Cells.Item(1, 8).Value  = "3일상승확률(%)"
$ws.Cells.Item(1, 9).Value  = "5일상승확률(%)"
$ws.Cells.Item(1, 10).Value = "10일상승확률(%)"

$ws.Cells.Item(1, 11).Value = "최종점수"
$ws.Cells.Item(1, 12).Value = "예측방식"
$ws.Cells.Item(1, 13).Value = "판단"
$ws.Cells.Item(1, 14).Value = "MACRO_SCORE"
$ws.Cells.Item(1, 15).Value = "MACRO_SIGNAL"

# Give the 4 brand-new header cells (L1:O1) the same look as the existing
# bold/bordered/centered header style (copy format only, from K1).
$ws.Cells.Item(1, 11).Copy()
$ws.Range($ws.Cells.Item(1, 12), $ws.Cells.Item(1, 15)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Data rows (rows 2-5): new tickers/order + new metric columns
#    Columns: A Date, B Name, C Ticker, D Close, E (cleared), F 5d-return,
#             G Score(rule), H 3d-prob, I 5d-prob, J 10d-prob, K FinalScore,
#             L PredictionMethod, M Judgement, N MacroScore, O MacroSignal
# ---------------------------------------------------------------------------
$rows = @(
    @{ B = "D-Wave Quantum Inc.";              C = "QBTS"; D = 22.67;  F = 10.56; G = 20; H = 60; I = 76; J = 73; K = 62;   M = "📈 매수 관찰 구간입니다." },
    @{ B = "International Business Machines";  C = "IBM";  D = 308.39; F = 6.19;  G = 40; H = 63; I = 60; J = 50; K = 61.6; M = "📈 매수 관찰 구간입니다." },
    @{ B = "Rigetti Computing, Inc.";          C = "RGTI"; D = 25.6;   F = 12.28; G = 20; H = 60; I = 70; J = 76; K = 59.6; M = "⛔ 관망하십시오." },
    @{ B = "IonQ, Inc.";                       C = "IONQ"; D = 49.22;  F = 20.05; G = 30; H = 56; I = 56; J = 70; K = 57;   M = "⛔ 관망하십시오." }
)

$newDate = "2025-11-29"
$r = 2
foreach ($row in $rows) {
    # --- A: date, forced to stay as literal text (not an Excel date serial) ---
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $newDate
    # strip the now-unneeded "@" text format back off by pasting the format
    # of a neighbouring plain General-formatted cell over it
    $ws.Cells.Item($r, 3).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).ClearContents()
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = "Pattern"
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = 85.36763896678245
    $ws.Cells.Item($r, 15).Value = "🟢 완화적 (상승 우위)"

    $r++
}
